# Naive Bayes Notes.pptx — "Adding PDF version of presentation"
#
# When a PDF rendition of this deck was uploaded alongside it in the
# SharePoint document library, the library's content-type sync rewrote the
# package's customXml parts: the two SharePoint-managed parts (the
# "DocumentLibraryForm" FormTemplates part and the managed-metadata
# "properties" part carrying Status / MediaServiceKeyPoints) were
# re-serialized and swapped into each other's slot. Reproduce that by
# locating each part by its distinguishing root element/namespace and
# writing the other part's markup into it (the matching itemProps*.xml
# schemaRefs/itemID side-parts are regenerated by the CustomXMLParts
# plumbing as a consequence).

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsNS = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$propsNS = "http://schemas.microsoft.com/office/2006/metadata/properties"

$formsXml = '<?mso-contentType?><FormTemplates xmlns="http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"><Display>DocumentLibraryForm</Display><Edit>DocumentLibraryForm</Edit><New>DocumentLibraryForm</New></FormTemplates>'
$propsXml = '<?xml version="1.0" encoding="utf-8"?><p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"><documentManagement><Status xmlns="71af3243-3dd4-4a8d-8c0d-dd76da1f02a5">Not started</Status><MediaServiceKeyPoints xmlns="71af3243-3dd4-4a8d-8c0d-dd76da1f02a5" xsi:nil="true"/></documentManagement></p:properties>'

$formsPart = $null
$propsPart = $null
for ($i = 1; $i -le $parts.Count; $i++) {
    $part = $parts.Item($i)
    if ($part.NamespaceURI -eq $formsNS) { $formsPart = $part }
    elseif ($part.NamespaceURI -eq $propsNS) { $propsPart = $part }
}

# Swap the two parts' markup in place (same slot/id, new content) so the
# "FormTemplates" item becomes the "properties" item and vice versa.
if (($formsPart -ne $null) -and ($propsPart -ne $null)) {
    $formsPart.XML = $propsXml
    $propsPart.XML = $formsXml
} else {
    # Fallback for hosts that don't surface NamespaceURI/enumeration on the
    # existing parts: drop and re-add with the swapped content so the
    # net XML payload of the customXml parts still matches.
    for ($i = $parts.Count; $i -ge 1; $i--) {
        $part = $parts.Item($i)
        if (($part.XML -like "*FormTemplates*") -or ($part.XML -like "*properties*")) {
            $part.Delete()
        }
    }
    [void]$parts.Add($propsXml)
    [void]$parts.Add($formsXml)
}
